$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 25 data
$ws.Range("A25").Value = "轮播图"
$ws.Range("B25").Value = 42806

# Copy style from A24/B24 to A25/B25
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B24").Copy()
$ws.Range("B25").PasteSpecial(-4122) # xlPasteFormats

# Set row height for row 25
$ws.Rows.Item(25).RowHeight = 35.25

# Update view: scroll to A19 and select E23
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("E23").Select()
